$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 16

$ws.Cells.Item($row, 1).Value = 11
$ws.Cells.Item($row, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item($row, 3).Value = "Bíobío"
$ws.Cells.Item($row, 4).Value = 45021
$ws.Cells.Item($row, 4).Style = $ws.Cells.Item(15, 4).Style
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item(15, 4).NumberFormat
$ws.Cells.Item($row, 5).Value = 8
$ws.Cells.Item($row, 6).Value = 100112052
$ws.Cells.Item($row, 7).Value = "Albahaca"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 50
$ws.Cells.Item($row, 11).Value = 4500
$ws.Cells.Item($row, 12).Value = 5000
$ws.Cells.Item($row, 13).Value = 4700
$ws.Cells.Item($row, 14).Value = "$/docena de matas"
$ws.Cells.Item($row, 15).Value = "Región Metropolitana"
$ws.Cells.Item($row, 16).Value = 783
$ws.Cells.Item($row, 17).Value = 6
$ws.Cells.Item($row, 18).Value = "Hortaliza"
